$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "303.51") are preserved as text, matching the source data which
# stores all prices/volumes as inline strings rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.933.34"
$ws.Range("D3").Value = "2.278.90"
$ws.Range("E3").Value = "  +4.13%  "
$ws.Range("D5").Value = "303.51"
$ws.Range("E5").Value = "  +4.37%  "
$ws.Range("D6").Value = "93.39"
$ws.Range("E6").Value = "  +8.34%  "
$ws.Range("D7").Value = "0.526"
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  +6.14%  "
$ws.Range("D10").Value = "32.95"
$ws.Range("E10").Value = "  +10.00%  "
$ws.Range("D11").Value = "54.68"
$ws.Range("E11").Value = "  +9.68%  "
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").Value = "6.71"
$ws.Range("E14").Value = "  +5.16%  "
$ws.Range("D15").Value = "2.631.66"
$ws.Range("E15").Value = "  +4.07%  "
$ws.Range("D16").Value = "14.28"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").Value = "2.270.69"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("E18").Value = "  +4.98%  "
$ws.Range("D19").Value = "41.831.60"
$ws.Range("E19").Value = "  +5.48%  "
$ws.Range("D20").Value = "12.37"
$ws.Range("E20").Value = "  +11.80%  "
$ws.Range("D21").Value = "0.0₃0913"
$ws.Range("E21").Value = "  +4.09%  "
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("D23").Value = "67.38"
$ws.Range("E23").Value = "  +3.71%  "
$ws.Range("D24").Value = "241.45"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("E25").Value = "  +7.21%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "1.89"
$ws.Range("E27").Value = "  +5.55%  "
$ws.Range("D28").Value = "23.95"
$ws.Range("E28").Value = "  +4.35%  "
$ws.Range("E29").Value = "  +7.70%  "
$ws.Range("D30").Value = "9.74"
$ws.Range("E30").Value = "  +6.85%  "
$ws.Range("E31").Value = "  +10.69%  "
$ws.Range("D32").Value = "158.41"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "5.22"
$ws.Range("E34").Value = "  +7.20%  "
$ws.Range("E35").Value = "  +5.72%  "
$ws.Range("D36").Value = "3.08"
$ws.Range("E36").Value = "  +9.93%  "
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("E38").Value = "  +12.35%  "
$ws.Range("E39").Value = "  +7.49%  "
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("E41").Value = "  +8.39%  "
$ws.Range("D42").Value = "3.99"
$ws.Range("E42").Value = "  +7.78%  "
$ws.Range("D43").Value = "2.062.73"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D44").Value = "20.05"
$ws.Range("E44").Value = "  +17.21%  "
$ws.Range("E45").Value = "  +5.32%  "
$ws.Range("D46").Value = "3.00"
$ws.Range("E46").Value = "  +13.54%  "
$ws.Range("D47").Value = "10.11"
$ws.Range("E47").Value = "  +4.54%  "
$ws.Range("E48").Value = "  -4.78%  "
$ws.Range("D49").Value = "2.499.65"
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("D50").Value = "1.52"
$ws.Range("E50").Value = "  +4.98%  "
$ws.Range("E51").Value = "  +5.84%  "
